# Update math practice problems in the division worksheet table.
$d = $word.ActiveDocument

$replacements = @(
    @("687÷8=", "702÷8="),
    @("305÷5=", "266÷2="),
    @("244÷2=", "588÷3="),
    @("930÷3=", "385÷2="),
    @("846÷5=", "987÷6="),
    @("134÷3=", "641÷9="),
    @("944÷2=", "383÷9="),
    @("708÷3=", "196÷9="),
    @("245÷9=", "498÷4="),
    @("502÷4=", "950÷7="),
    @("923÷2=", "922÷3="),
    @("606÷9=", "415÷7="),
    @("768÷6=", "931÷5="),
    @("160÷7=", "163÷4="),
    @("109÷4=", "180÷5="),
    @("999÷9=", "135÷2="),
    @("755÷7=", "205÷4="),
    @("740÷8=", "233÷6="),
    @("281÷9=", "205÷9="),
    @("183÷8=", "195÷9="),
    @("540÷5=", "592÷3="),
    @("898÷2=", "921÷6="),
    @("768÷9=", "907÷4="),
    @("857÷6=", "797÷6="),
    @("308÷8=", "791÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
